# TC_71808 - Battery Standby for FC test cases update
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels")

# --- Move "AlarmLoadingDetail" / "StandbyLoadingDetail" header labels ---
# from S7:T7 (old extra columns) up to F1:G1
$ws.Range("S7:T7").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "AlarmLoadingDetail"
$ws.Range("G1").Value = "StandbyLoadingDetail"

# --- Move "Battery Alarm (A)" / "Battery Standby (A)" values ---
# from S8:T8 (old extra columns) up to F2:G2
$ws.Range("S8:T8").Copy()
$ws.Range("F2").PasteSpecial(-4122)
$ws.Range("F2").Value = "Battery Alarm (A)"
$ws.Range("G2").Value = "Battery Standby (A)"

# Resize column G (7) to the width that column T (20) used to have,
# since column T's content moved into column G
$ws.Columns.Item(7).ColumnWidth = $ws.Columns.Item(20).ColumnWidth

# Clear the now unused columns S:T (rows 7 and 8) that held this data
$ws.Range("S7:T8").Clear()

# Update selection / view to reflect the now-empty S:T columns
$ws.Range("S1:T1048576").Select()

$excel.CutCopyMode = $false
